# Updated cryptos list on Sun Sep 17 14:17:59 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# the coinranking.com snapshot, and rewrites the bottom three rows because
# a new coin (BabyDogeCoin) was inserted into the ranking, pushing
# EnergySwap / Algorand down a row and dropping Mantle off the bottom.
#
# Numeric-looking price strings must stay text (the sheet stores "Price" as
# plain text, e.g. "26.756.26", "217.19", "2.37"), so NumberFormat is forced
# to "@" immediately before assigning those values to stop Excel from
# re-interpreting them as numbers/dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Cells.Item(2, 4).Value = '26.756.26'
$ws.Cells.Item(2, 5).Value = '  +0.28%  '

# Row 3 - Ethereum
$ws.Cells.Item(3, 4).Value = '1.641.24'
$ws.Cells.Item(3, 5).Value = '  -0.09%  '

# Row 4 - TetherUSD
$ws.Cells.Item(4, 5).Value = '  +0.21%  '

# Row 5 - BNB
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '217.19'
$ws.Cells.Item(5, 5).Value = '  +1.02%  '

# Row 6 - XRP
$ws.Cells.Item(6, 5).Value = '  -0.24%  '

# Row 7 - USDC
$ws.Cells.Item(7, 5).Value = '  +0.26%  '

# Row 8 - Cardano
$ws.Cells.Item(8, 5).Value = '  -0.19%  '

# Row 9 - Dogecoin
$ws.Cells.Item(9, 5).Value = '  +0.01%  '

# Row 10 - Solana
$ws.Cells.Item(10, 5).Value = '  +0.46%  '

# Row 11 - TRON
$ws.Cells.Item(11, 5).Value = '  -0.03%  '

# Row 12 - WrappedliquidstakedEther2.0
$ws.Cells.Item(12, 4).Value = '1.870.33'
$ws.Cells.Item(12, 5).Value = '  -0.06%  '

# Row 13 - WrappedEther
$ws.Cells.Item(13, 4).Value = '1.648.62'
$ws.Cells.Item(13, 5).Value = '  +1.50%  '

# Row 14 - Polkadot
$ws.Cells.Item(14, 5).Value = '  -0.63%  '

# Row 15 - Polygon
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '0.527'
$ws.Cells.Item(15, 5).Value = '  -0.38%  '

# Row 16 - Litecoin
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '64.49'
$ws.Cells.Item(16, 5).Value = '  -0.92%  '

# Row 17 - WrappedBTC
$ws.Cells.Item(17, 4).Value = '26.737.75'
$ws.Cells.Item(17, 5).Value = '  +0.19%  '

# Row 18 - ShibaInu
$ws.Cells.Item(18, 5).Value = '  -1.38%  '

# Row 19 - BitcoinCash
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '213.83'
$ws.Cells.Item(19, 5).Value = '  -1.17%  '

# Row 20 - Dai
$ws.Cells.Item(20, 5).Value = '  +0.29%  '

# Row 21 - Uniswap
$ws.Cells.Item(21, 5).Value = '  +0.65%  '

# Row 22 - Toncoin
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '2.37'
$ws.Cells.Item(22, 5).Value = '  +5.02%  '

# Row 23 - Chainlink
$ws.Cells.Item(23, 5).Value = '  -0.76%  '

# Row 24 - Avalanche
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '9.27'
$ws.Cells.Item(24, 5).Value = '  -2.33%  '

# Row 25
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '145.51'
$ws.Cells.Item(25, 5).Value = '  +0.08%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.03%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -1.21%  '

# Row 28 - Cosmos
$ws.Cells.Item(28, 5).Value = '  +0.16%  '

# Row 29 - EthereumClassic
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '15.60'
$ws.Cells.Item(29, 5).Value = '  -0.73%  '

# Row 30 - Hedera
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0507'
$ws.Cells.Item(30, 5).Value = '  -1.41%  '

# Row 31 - PancakeSwap
$ws.Cells.Item(31, 5).Value = '  +0.76%  '

# Row 33 - InternetComputer(DFINITY)
$ws.Cells.Item(33, 5).Value = '  -1.08%  '

# Row 34 - Maker
$ws.Cells.Item(34, 4).Value = '1.285.73'
$ws.Cells.Item(34, 5).Value = '  +0.47%  '

# Row 35 - LidoDAOToken
$ws.Cells.Item(35, 5).Value = '  -0.47%  '

# Row 36 - HuobiToken
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '2.44'
$ws.Cells.Item(36, 5).Value = '  +1.47%  '

# Row 37 - VeChain
$ws.Cells.Item(37, 5).Value = '  -0.84%  '

# Row 38 - ImmutableX
$ws.Cells.Item(38, 5).Value = '  +0.60%  '

# Row 39 - ARBITRUM
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.817'
$ws.Cells.Item(39, 5).Value = '  -1.53%  '

# Row 40 - PaxDollar
$ws.Cells.Item(40, 5).Value = '  +0.32%  '

# Row 41 - TrustWalletToken
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.805'
$ws.Cells.Item(41, 5).Value = '  -1.02%  '

# Row 42 - MXToken
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '2.22'
$ws.Cells.Item(42, 5).Value = '  -1.11%  '

# Row 43 - FraxShare
$ws.Cells.Item(43, 5).Value = '  -2.62%  '

# Row 44 - RocketPoolETH
$ws.Cells.Item(44, 4).Value = '1.781.44'
$ws.Cells.Item(44, 5).Value = '  -0.01%  '

# Row 45 - Aave
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '61.40'
$ws.Cells.Item(45, 5).Value = '  +3.49%  '

# Row 46 - Quant
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '91.84'

# Row 47 - RenderToken
$ws.Cells.Item(47, 5).Value = '  +0.11%  '

# Row 48 - Cronos
$ws.Cells.Item(48, 5).Value = '  +0.32%  '

# Row 49 - new entrant BabyDogeCoin (was EnergySwap)
$ws.Cells.Item(49, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(49, 4).Value = '0.0₆0101'
$ws.Cells.Item(49, 5).Value = '  -3.97%  '

# Row 50 - EnergySwap (was Algorand)
$ws.Cells.Item(50, 2).Value = 'EnergySwap'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.67'
$ws.Cells.Item(50, 5).Value = '  -1.39%  '

# Row 51 - Algorand (was Mantle)
$ws.Cells.Item(51, 2).Value = 'Algorand'
$ws.Cells.Item(51, 3).Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.0966'
$ws.Cells.Item(51, 5).Value = '  +0.01%  '
